# Match the width of the logo in the header to be the minimum standard.
#
# The primary (default) header contains a single inline picture (the
# company logo). Resize it from 648000 x 339497 EMU (~51.02pt x 26.73pt)
# up to 1274400 x 667677 EMU (~100.35pt x 52.57pt), preserving its
# aspect ratio.

$d = $word.ActiveDocument

$targetWidthPt  = 1274400 / 914400 * 72   # 100.34645669291339
$targetHeightPt = 667677  / 914400 * 72   # 52.572992125984255

$section = $d.Sections.First
$header  = $section.Headers.Item(1)       # wdHeaderFooterPrimary

$logo = $header.Range.InlineShapes.Item(1)
$logo.LockAspectRatio = $true
$logo.Width  = $targetWidthPt
$logo.Height = $targetHeightPt

Write-Output ("Logo resized to " + $logo.Width + "pt x " + $logo.Height + "pt")
